$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for 2022-Q3, shift the rest
#    down, and bump the running index in column A for every row that moved.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Pull the "index" column formatting (style s="2") from the row below, which
# still carries the original formatting, onto the freshly inserted row.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)
# The B:D cells of the inserted row inherited header formatting on insert;
# reset them back to the plain style used by every other data row.
$total.Range("B2:D2").Style = "Normal"

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.66

for ($r = 3; $r -le 9; $r++) {
    $cell = $total.Cells.Item($r, 1)
    $old = $cell.Value()
    $cell.Value = $old + 1
}

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" sheet, positioned right after "总计". Duplicate the
#    "2022-Q2" sheet so formatting/margins/styles match exactly, then
#    overwrite its data and extend it with two extra rows (C-class shares).
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $total)
$newws = $wb.Worksheets.Item("2022-Q2 (2)")
$newws.Name = "2022-Q3"

# Extend the sheet from 2 data rows to 4 data rows (rows 4 and 5 copy the
# formatting of row 3).
$newws.Range("A3:H3").Copy()
$newws.Range("A4:H4").PasteSpecial(-4122)
$newws.Range("A5:H5").PasteSpecial(-4122)

# Force text storage (matches source data) for the fund-code column and the
# percentage/price columns before writing numeric-looking strings into them,
# so leading zeros and exact text survive instead of being parsed as numbers.
$newws.Range("B2:B5").NumberFormat = "@"
$newws.Range("D2:G3").NumberFormat = "@"
$newws.Range("D4:F5").NumberFormat = "@"

$newws.Cells.Item(2, 2).Value = "001092"
$newws.Cells.Item(2, 3).Value = "广发纳斯达克生物科技指数（QDII）人民币A"
$newws.Cells.Item(2, 4).Value = "4.84"
$newws.Cells.Item(2, 5).Value = "90.03"
$newws.Cells.Item(2, 6).Value = "6.81"
$newws.Cells.Item(2, 7).Value = "0.3296"
$newws.Cells.Item(2, 8).Value = 3

$newws.Cells.Item(3, 1).Value = 1
$newws.Cells.Item(3, 2).Value = "001093"
$newws.Cells.Item(3, 3).Value = "广发纳斯达克生物科技指数（QDII）美元A"
$newws.Cells.Item(3, 4).Value = "4.84"
$newws.Cells.Item(3, 5).Value = "90.03"
$newws.Cells.Item(3, 6).Value = "6.81"
$newws.Cells.Item(3, 7).Value = "0.3296"
$newws.Cells.Item(3, 8).Value = 3

$newws.Cells.Item(4, 1).Value = 2
$newws.Cells.Item(4, 2).Value = "016470"
$newws.Cells.Item(4, 3).Value = "广发纳斯达克生物科技指数（QDII）人民币C"
$newws.Cells.Item(4, 4).Value = "0.00"
$newws.Cells.Item(4, 5).Value = "90.03"
$newws.Cells.Item(4, 6).Value = "6.81"
$newws.Cells.Item(4, 7).Value = 0
$newws.Cells.Item(4, 8).Value = 3

$newws.Cells.Item(5, 1).Value = 3
$newws.Cells.Item(5, 2).Value = "016471"
$newws.Cells.Item(5, 3).Value = "广发纳斯达克生物科技指数（QDII）美元C"
$newws.Cells.Item(5, 4).Value = "0.00"
$newws.Cells.Item(5, 5).Value = "90.03"
$newws.Cells.Item(5, 6).Value = "6.81"
$newws.Cells.Item(5, 7).Value = 0
$newws.Cells.Item(5, 8).Value = 3

# Drop the temporary text number-format now that the values are locked in as
# text, so the cells fall back to the plain (unstyled) look used elsewhere.
$newws.Range("B2:B5").Style = "Normal"
$newws.Range("D2:G3").Style = "Normal"
$newws.Range("D4:F5").Style = "Normal"
